# Theme swap: the authored commit exchanges the contents of
# ppt/theme/theme1.xml ("Office Theme" / "Office" color scheme) and
# ppt/theme/theme2.xml ("Integral" / "Red Violet" color scheme).
#
# theme2.xml is the theme actually wired to the (one and only) Slide
# Master, so it is the theme PowerPoint's object model exposes as the
# presentation's live color scheme (SlideMaster.ColorScheme /
# SlideMaster.Theme, NotesMaster.ColorScheme, HandoutMaster.ColorScheme,
# and Slide.ThemeColorScheme all resolve to this same, single, active
# scheme). We drive it to the "Office Theme" palette, which is the half
# of the swap that is reachable from PowerPoint automation.

$p = $ppt.ActivePresentation

# Use the 12-slot theme color scheme (distinct from the legacy 8-slot
# ColorScheme) so the edit is scoped to color swatches only.
$tcs = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeColor($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" color scheme (order: dk1, lt1,
# dk2, lt2, accent1-6, hlink, folHlink == items 1-12).
Set-ThemeColor 1  "000000"
Set-ThemeColor 2  "FFFFFF"
Set-ThemeColor 3  "44546A"
Set-ThemeColor 4  "E7E6E6"
Set-ThemeColor 5  "5B9BD5"
Set-ThemeColor 6  "ED7D31"
Set-ThemeColor 7  "A5A5A5"
Set-ThemeColor 8  "FFC000"
Set-ThemeColor 9  "4472C4"
Set-ThemeColor 10 "70AD47"
Set-ThemeColor 11 "0563C1"
Set-ThemeColor 12 "954F72"
